# DaySale report: a new shortage item (ZURCAL) was added to the list.
# That pushes the previous 5th item ("اولويز ماكس طويل جدا") down to a new
# 6th row, shifts the totals row and the footer row down by one, and bumps
# the running total (P column) by the new item's price. The report's
# generated timestamp (footer, row A) is also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert a blank row before the current totals row (12).
#        This pushes: old row 12 (totals)  -> row 13
#                      old row 13 (footer)  -> row 14
#        Row 11 (the last item row, "اولويز ماكس طويل جدا") is untouched.
$ws.Rows(12).Insert()

# --- 2. Populate the new row 12 by cloning row 11's layout (same column
#        merges / cell styles as every other item row), then restore the
#        two row heights the diff calls for.
$ws.Range("A11:Q11").Copy($ws.Range("A12:Q12"))
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 24.75

# Row 12 now holds an exact copy of row 11 ("اولويز ماكس طويل جدا" with its
# numbers) - that is exactly what should live there, it just needs the
# sequence number bumped from 5 to 6.
$ws.Range("A12").Value2 = 6

# --- 3. Turn row 11 into the new ZURCAL item.
#        Re-use existing identical text cells where possible (copy) so no
#        stray formatting/styles are introduced; only genuinely new text
#        is typed in directly.
$ws.Range("C11").Value2 = "ZURCAL 40MG 14 GASTRO RESISTANT TAB"
$ws.Range("H9").Copy($ws.Range("H11"))     # "4:0" (already used by H9)
$ws.Range("L8").Copy($ws.Range("L11"))     # "1"   (already used by L8/L9/L10)
$ws.Range("N11").Value2 = "96.00"

# P11 ("96.0000") needs to stay literal text (like every other price cell
# in this column) even though the cell's number format is numeric - force
# text entry, then restore the original numeric-style formatting so the
# style index matches the rest of the column.
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value2 = "96.0000"
$ws.Range("P10").Copy()
$ws.Range("P11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Q11 ("1:0") is already correct and untouched by the insert.

# --- 4. Update the totals row (now row 13): add the new item's price.
$ws.Range("P13").Value2 = 256.045

# --- 5. Refresh the footer timestamp (now row 14).
$ws.Range("A14").Value2 = "Monday, 25 August, 2025 10:22 AM"

Write-Output "done"
